# Added script for end to end testing:
# create a new "testdata" worksheet (after the existing "logindata" sheet)
# containing Firstname / Lastname / zipcode sample data.

$wb = $excel.ActiveWorkbook
$loginSheet = $wb.Worksheets.Item("logindata")

# Insert the new worksheet right after "logindata" so it becomes the
# second (and active) sheet in the workbook.
$ws = $wb.Worksheets.Add($null, $loginSheet)
$ws.Name = "testdata"

$ws.Range("A1").Value = "Firstname"
$ws.Range("B1").Value = "Lastname"
$ws.Range("C1").Value = "zipcode"

$ws.Range("A2").Value = "pintu"
$ws.Range("B2").Value = "samal"
$ws.Range("C2").Value = 755004

# Make the newly added sheet the active one with C2 selected, matching
# the state the workbook was saved in.
$ws.Range("C2").Select()
